$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string / label text updates ---
$ws.Range("B34").Value = "Soudan du Sud*"   # resource-rich flag added
$ws.Range("B48").Value = "Cabo Verde*"      # resource-rich flag added
$ws.Range("B57").Value = "Nigeria"          # resource-rich flag removed

# --- Row formatting updates (resource-rich shading follows the label change) ---
# Row 34 (Soudan du Sud) gains the resource-rich shading used by other "*" rows (e.g. row 17, Tchad*)
$ws.Range("B17:AE17").Copy()
$ws.Range("B34:AE34").PasteSpecial(-4122)
# Row 57 (Nigeria) loses the resource-rich shading, matching a plain row (e.g. row 5, Botswana)
$ws.Range("B5:AE5").Copy()
$ws.Range("B57:AE57").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Recalculated aggregate values (regional/group totals affected by the Nigeria reclassification) ---
# Row 69
$ws.Range("C69").Value = 481513.3
$ws.Range("D69").Value = 460479
$ws.Range("E69").Value = 437.11849999999998
$ws.Range("F69").Value = 147233.9
$ws.Range("G69").Value = 51752.6
$ws.Range("H69").Value = 44425.9
$ws.Range("I69").Value = 185561.87
$ws.Range("J69").Value = 20943.9512
$ws.Range("K69").Value = 662321.16449999996
$ws.Range("L69").Value = 247767.9442
$ws.Range("M69").Value = 697847.64190000005
$ws.Range("O69").Value = 30530.7945
$ws.Range("P69").Value = 1301.6143
$ws.Range("Q69").Value = 1843.0345
$ws.Range("W69").Value = 756.98472181198099
$ws.Range("X69").Value = 215447.17600000001
$ws.Range("Y69").Value = 8098.42
$ws.Range("AC69").Value = 404.21643
$ws.Range("AD69").Value = 231.4127
$ws.Range("AE69").Value = 64.170810000000003

# Row 77
$ws.Range("C77").Value = 1679542.105
$ws.Range("D77").Value = 1648469.2
$ws.Range("E77").Value = 3583.8971000000001
$ws.Range("F77").Value = 509096.4632
$ws.Range("G77").Value = 128012.9942
$ws.Range("H77").Value = 110868.32859999999
$ws.Range("I77").Value = 789952.06
$ws.Range("J77").Value = 89123.712799999994
$ws.Range("K77").Value = 680091.77020000003
$ws.Range("L77").Value = 936381.03350000002
$ws.Range("M77").Value = 1063091.0549999999
$ws.Range("N77").Value = -363267.0968
$ws.Range("O77").Value = 22444.962899999999
$ws.Range("P77").Value = 41573.5726
$ws.Range("Q77").Value = 50675.8675
$ws.Range("W77").Value = 957.74728145051904
$ws.Range("X77").Value = 497142.42300000001
$ws.Range("Y77").Value = 63939.044000000002
$ws.Range("Z77").Value = 1109.0487000000001
$ws.Range("AA77").Value = 3724.0740000000001
$ws.Range("AB77").Value = 1266.4834900000001
$ws.Range("AC77").Value = 8997.2543999999998
$ws.Range("AD77").Value = 9070.8014000000003
$ws.Range("AE77").Value = 8365.4426000000003

# Row 80
$ws.Range("C80").Value = 899131.4
$ws.Range("D80").Value = 894106.1
$ws.Range("E80").Value = 557.13019999999995
$ws.Range("F80").Value = 244919.497
$ws.Range("G80").Value = 25269.316999999999
$ws.Range("H80").Value = 23183.316999999999
$ws.Range("I80").Value = 128526.92
$ws.Range("J80").Value = 7739.2939999999999
$ws.Range("K80").Value = 75866.051200000002
$ws.Range("L80").Value = 200041.14360000001
$ws.Range("M80").Value = 90340.597599999994
$ws.Range("N80").Value = -157.85919999999999
$ws.Range("O80").Value = 14697.933199999999
$ws.Range("P80").Value = 632.02729999999997
$ws.Range("Q80").Value = 1988.9505999999999
$ws.Range("R80").Value = 31141.521000000001
$ws.Range("S80").Value = 31127.921999999999
$ws.Range("T80").Value = 22965.463
$ws.Range("U80").Value = 8175.9949999999999
$ws.Range("V80").Value = 13.534000000000001
$ws.Range("W80").Value = 704.92833769968695
$ws.Range("X80").Value = 43997.116000000002
$ws.Range("Y80").Value = 2545.808
$ws.Range("Z80").Value = 1164.2
$ws.Range("AC80").Value = 114.95764
$ws.Range("AD80").Value = 79.073539999999994
$ws.Range("AE80").Value = 41.683369999999996

# Row 82
$ws.Range("C82").Value = 2111307.5871000001
$ws.Range("D82").Value = 2071140.6094
$ws.Range("E82").Value = 3917.4252999999999
$ws.Range("F82").Value = 873956.16200000001
$ws.Range("G82").Value = 256527.39
$ws.Range("H82").Value = 222346.69
$ws.Range("I82").Value = 507332.65
$ws.Range("J82").Value = 42580.103199999998
$ws.Range("K82").Value = 1107886.1105
$ws.Range("L82").Value = 938579.94880000001
$ws.Range("M82").Value = 1129152.112
$ws.Range("N82").Value = 43308.144099999998
$ws.Range("O82").Value = 86040.679099999994
$ws.Range("P82").Value = 22117.452000000001
$ws.Range("Q82").Value = 30752.5802
$ws.Range("R82").Value = 253095.571
$ws.Range("S82").Value = 244294.99600000001
$ws.Range("T82").Value = 196108.997
$ws.Range("U82").Value = 54113.608999999997
$ws.Range("V82").Value = 5932.6149999999998
$ws.Range("W82").Value = 715.29824994873502
$ws.Range("X82").Value = 766322.728
$ws.Range("Y82").Value = 40452.120999999999
$ws.Range("Z82").Value = 6695.7455600000003
$ws.Range("AC82").Value = 4262.5187400000004
$ws.Range("AD82").Value = 1651.7487000000001
$ws.Range("AE82").Value = 942.62345000000005

# Row 84
$ws.Range("C84").Value = 1454089.7871000001
$ws.Range("D84").Value = 1425017.2094000001
$ws.Range("E84").Value = 1120.288
$ws.Range("F84").Value = 532378.62199999997
$ws.Range("G84").Value = 136292.34
$ws.Range("H84").Value = 125838.34
$ws.Range("I84").Value = 341879.76
$ws.Range("J84").Value = 31254.030599999998
$ws.Range("K84").Value = 851495.00890000002
$ws.Range("L84").Value = 676581.09510000004
$ws.Range("M84").Value = 902881.68900000001
$ws.Range("O84").Value = 78770.378899999996
$ws.Range("P84").Value = 6034.6152000000002
$ws.Range("Q84").Value = 4896.6900999999998
$ws.Range("R84").Value = 72317.024999999994
$ws.Range("S84").Value = 70638.099000000002
$ws.Range("T84").Value = 55392.847000000002
$ws.Range("U84").Value = 16030.477999999999
$ws.Range("W84").Value = 639.38306112763303
$ws.Range("X84").Value = 468265.842
$ws.Range("Y84").Value = 19744.322
$ws.Range("AC84").Value = 1153.6528599999999
$ws.Range("AD84").Value = 405.59609
$ws.Range("AE84").Value = 175.82222999999999

# Row 86
$ws.Range("C86").Value = 1088062.2
$ws.Range("D86").Value = 1075143.5
$ws.Range("E86").Value = 2341.5455000000002
$ws.Range("F86").Value = 407645.32699999999
$ws.Range("G86").Value = 129210.217
$ws.Range("H86").Value = 104372.917
$ws.Range("I86").Value = 228767.63
$ws.Range("J86").Value = 13139.759899999999
$ws.Range("K86").Value = 265549.31050000002
$ws.Range("L86").Value = 407918.82789999997
$ws.Range("M86").Value = 249328.54680000001
$ws.Range("O86").Value = 21715.793300000001
$ws.Range("P86").Value = 13998.4673
$ws.Range("Q86").Value = 17452.649300000001
$ws.Range("R86").Value = 190188.78400000001
$ws.Range("S86").Value = 183831.12299999999
$ws.Range("T86").Value = 151697.71599999999
$ws.Range("U86").Value = 36573.720999999998
$ws.Range("W86").Value = 792.59854707928002
$ws.Range("X86").Value = 296055.89500000002
$ws.Range("Y86").Value = 17495.900000000001
$ws.Range("AC86").Value = 2806.6102900000001
$ws.Range("AD86").Value = 1061.1628000000001
$ws.Range("AE86").Value = 667.72004000000004

# Row 87
$ws.Range("C87").Value = 1475168.6518999999
$ws.Range("D87").Value = 1408358.4802999999
$ws.Range("E87").Value = 7888.8406999999997
$ws.Range("F87").Value = 635275.10329999996
$ws.Range("G87").Value = 378571.6263
$ws.Range("H87").Value = 318190.14630000002
$ws.Range("I87").Value = 397430.02
$ws.Range("J87").Value = 29860.411499999998
$ws.Range("K87").Value = 415114.94520000002
$ws.Range("L87").Value = 2187058.6205000002
$ws.Range("M87").Value = 921455.05110000004
$ws.Range("N87").Value = -48997.489500000003
$ws.Range("O87").Value = 28930.269100000001
$ws.Range("P87").Value = 168352.88430000001
$ws.Range("Q87").Value = 44418.7264
$ws.Range("R87").Value = 836346.73600000003
$ws.Range("S87").Value = 787159.72199999995
$ws.Range("T87").Value = 586205.16799999995
$ws.Range("U87").Value = 217154.231
$ws.Range("V87").Value = 17202.02
$ws.Range("W87").Value = 503.70847914792199
$ws.Range("X87").Value = 825820.12500000105
$ws.Range("Y87").Value = 71513.097999999998
$ws.Range("AA87").Value = 6877.2483400000001
$ws.Range("AC87").Value = 36008.686229999999
$ws.Range("AD87").Value = 14287.245870000001
$ws.Range("AE87").Value = 7127.9696400000003

# Row 89
$ws.Range("C89").Value = 4991019.9950000001
$ws.Range("D89").Value = 4862576.55
$ws.Range("E89").Value = 21364.522000000001
$ws.Range("F89").Value = 1658452.6081999999
$ws.Range("G89").Value = 520006.26919999998
$ws.Range("H89").Value = 461773.79859999998
$ws.Range("I89").Value = 1929566.92
$ws.Range("J89").Value = 142418.04459999999
$ws.Range("K89").Value = -629108.59069999994
$ws.Range("L89").Value = 2305141.6275999998
$ws.Range("M89").Value = 1109561.5895
$ws.Range("N89").Value = -1717869.1196000001
$ws.Range("O89").Value = 23611.3727
$ws.Range("P89").Value = 250283.09270000001
$ws.Range("Q89").Value = 247824.9743
$ws.Range("R89").Value = 2112074.35
$ws.Range("S89").Value = 2012847.335
$ws.Range("T89").Value = 1403397.6629999999
$ws.Range("U89").Value = 660914.02500000002
$ws.Range("V89").Value = 53850.137999999999
$ws.Range("W89").Value = 784.68590995808199
$ws.Range("X89").Value = 1715551.3629999999
$ws.Range("Y89").Value = 315376.26500000001
$ws.Range("Z89").Value = 49777.824099999998
$ws.Range("AA89").Value = 23053.506509999999
$ws.Range("AC89").Value = 44265.617639999997
$ws.Range("AD89").Value = 21837.107759999999
$ws.Range("AE89").Value = 21525.686129999998

# Row 90
$ws.Range("C90").Value = 3700044.4685
$ws.Range("D90").Value = 3476299.7223999999
$ws.Range("E90").Value = 26002.997599999999
$ws.Range("F90").Value = 1224417.8367000001
$ws.Range("G90").Value = 359717.679
$ws.Range("H90").Value = 342872.41759999999
$ws.Range("I90").Value = 1030987.36
$ws.Range("J90").Value = 63906.263099999996
$ws.Range("K90").Value = -706217.06409999996
$ws.Range("L90").Value = 1498000.5881000001
$ws.Range("M90").Value = 204629.11569999999
$ws.Range("N90").Value = -910746.16799999995
$ws.Range("O90").Value = 24283.285400000001
$ws.Range("P90").Value = 157063.08929999999
$ws.Range("Q90").Value = 204023.00870000001
$ws.Range("R90").Value = 902114.94700000004
$ws.Range("S90").Value = 895624.12800000003
$ws.Range("T90").Value = 569637.41399999999
$ws.Range("U90").Value = 331026.451
$ws.Range("V90").Value = 7279.4440000000004
$ws.Range("W90").Value = 563.60999066376405
$ws.Range("X90").Value = 1895359.899
$ws.Range("Y90").Value = 499415.815
$ws.Range("Z90").Value = 38497.305800000002
$ws.Range("AA90").Value = 10128.014279999999
$ws.Range("AC90").Value = 28390.542870000001
$ws.Range("AD90").Value = 10136.66042
$ws.Range("AE90").Value = 9428.1796100000101

# Row 91
$ws.Range("O91").Value = 98667.441800000102

# Row 94
$ws.Range("C94").Value = 116517.8
$ws.Range("D94").Value = 111741.8
$ws.Range("E94").Value = 302.5401
$ws.Range("F94").Value = 15063.76
$ws.Range("G94").Value = 8859.33
$ws.Range("H94").Value = 6180.62
$ws.Range("I94").Value = 83283.89
$ws.Range("K94").Value = 74034.777499999997
$ws.Range("L94").Value = 80794.415399999998
$ws.Range("N94").Value = 37073.775500000003
$ws.Range("P94").Value = 1180.0676000000001
$ws.Range("Q94").Value = 8304.7536
$ws.Range("X94").Value = 24330.088
$ws.Range("Y94").Value = 1250.4459999999999
$ws.Range("Z94").Value = 756.33403999999996
$ws.Range("AC94").Value = 214.71969999999999
$ws.Range("AD94").Value = 61.313420000000001
$ws.Range("AE94").Value = 102.53394

# Row 97
$ws.Range("C97").Value = 2263070.8870999999
$ws.Range("D97").Value = 2226634.3094000001
$ws.Range("E97").Value = 2483.3953999999999
$ws.Range("F97").Value = 841053.18200000003
$ws.Range("G97").Value = 227087.44
$ws.Range("H97").Value = 200436.74
$ws.Range("I97").Value = 546204.42000000004
$ws.Range("J97").Value = 43701.464899999999
$ws.Range("K97").Value = 1140108.7282
$ws.Range("L97").Value = 974963.31759999995
$ws.Range("M97").Value = 1143686.2660000001
$ws.Range("N97").Value = 75005.770999999993
$ws.Range("O97").Value = 100145.6009
$ws.Range("P97").Value = 9991.1857999999993
$ws.Range("Q97").Value = 12398.6978
$ws.Range("R97").Value = 154222.285
$ws.Range("S97").Value = 147150.228
$ws.Range("T97").Value = 127411.351
$ws.Range("W97").Value = 679.58305034906698
$ws.Range("X97").Value = 655103.78300000005
$ws.Range("Y97").Value = 29752.843000000001
$ws.Range("AC97").Value = 2113.4783200000002
$ws.Range("AD97").Value = 869.43469000000005
$ws.Range("AE97").Value = 474.57981999999998

# Row 98
$ws.Range("C98").Value = 813249.2
$ws.Range("D98").Value = 785534.8
$ws.Range("E98").Value = 2201.0551999999998
$ws.Range("F98").Value = 278206.40000000002
$ws.Range("G98").Value = 114222.39999999999
$ws.Range("H98").Value = 101888.9
$ws.Range("I98").Value = 182426.96
$ws.Range("J98").Value = 15123.371999999999
$ws.Range("K98").Value = 258624.19940000001
$ws.Range("L98").Value = 700694.67689999996
$ws.Range("M98").Value = 230862.9859
$ws.Range("O98").Value = 11846.736199999999
$ws.Range("Q98").Value = 33296.177300000003
$ws.Range("R98").Value = 199175.66
$ws.Range("S98").Value = 194020.916
$ws.Range("T98").Value = 148937.41800000001
$ws.Range("U98").Value = 47609.074999999997
$ws.Range("V98").Value = 3117.0479999999998
$ws.Range("W98").Value = 469.73362772413498
$ws.Range("X98").Value = 203323.519
$ws.Range("Z98").Value = 6627.9411600000003
$ws.Range("AA98").Value = 748.98644000000002
$ws.Range("AB98").Value = 87.832999999999998
$ws.Range("AC98").Value = 7462.4008199999998
$ws.Range("AD98").Value = 2652.9184799999998
$ws.Range("AE98").Value = 985.11247000000003

